$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking price values stay as text (matches original inlineStr formatting)
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '29.365.99'
$ws.Range('E2').Value = '  -0.18%  '
$ws.Range('D3').Value = '1.847.69'
$ws.Range('E3').Value = '  -0.12%  '
$ws.Range('D4').Value = '0.9982'
$ws.Range('E4').Value = '  -0.17%  '
$ws.Range('D5').Value = '240.48'
$ws.Range('E5').Value = '  -0.12%  '
$ws.Range('D6').Value = '0.6271'
$ws.Range('E6').Value = '  -0.59%  '
$ws.Range('D7').Value = '0.9992'
$ws.Range('E7').Value = '  -0.13%  '
$ws.Range('D8').Value = '0.07625'
$ws.Range('E8').Value = '  -1.05%  '
$ws.Range('D9').Value = '0.2902'
$ws.Range('E9').Value = '  -1.38%  '
$ws.Range('D10').Value = '24.78'
$ws.Range('E10').Value = '  +1.00%  '
$ws.Range('D11').Value = '0.07735'
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '5.028'
$ws.Range('E12').Value = '  +0.00%  '
$ws.Range('D13').Value = '0.6795'
$ws.Range('E13').Value = '  -0.17%  '
$ws.Range('D14').Value = '0.00001053'
$ws.Range('E14').Value = '  -2.90%  '
$ws.Range('D15').Value = '83.02'
$ws.Range('E15').Value = '  -0.86%  '
$ws.Range('D16').Value = '6.144'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('D17').Value = '29.363.59'
$ws.Range('E17').Value = '  -0.28%  '
$ws.Range('D18').Value = '227.83'
$ws.Range('E18').Value = '  -0.65%  '
$ws.Range('D19').Value = '12.35'
$ws.Range('E19').Value = '  -0.99%  '
$ws.Range('E20').Value = '  -0.19%  '
$ws.Range('D21').Value = '7.462'
$ws.Range('E21').Value = '  +0.09%  '
$ws.Range('D22').Value = '0.9994'
$ws.Range('E22').Value = '  -0.13%  '
$ws.Range('D23').Value = '158.87'
$ws.Range('E23').Value = '  +0.86%  '
$ws.Range('D24').Value = '0.1386'
$ws.Range('E24').Value = '  -0.32%  '
$ws.Range('D25').Value = '8.438'
$ws.Range('E25').Value = '  +0.88%  '
$ws.Range('D26').Value = '17.66'
$ws.Range('E26').Value = '  -0.17%  '
$ws.Range('D27').Value = '1.408'
$ws.Range('E27').Value = '  +7.35%  '
$ws.Range('D28').Value = '1.461'
$ws.Range('E28').Value = '  -0.66%  '
$ws.Range('D29').Value = '0.05605'
$ws.Range('E29').Value = '  -2.26%  '
$ws.Range('D30').Value = '4.111'
$ws.Range('E30').Value = '  +0.02%  '
$ws.Range('D31').Value = '4.065'
$ws.Range('E31').Value = '  +0.27%  '
$ws.Range('E32').Value = '  +0.21%  '
$ws.Range('E33').Value = '  -1.10%  '
$ws.Range('D34').Value = '0.6962'
$ws.Range('E34').Value = '  -2.03%  '
$ws.Range('D35').Value = '2.588'
$ws.Range('E35').Value = '  +0.00%  '
$ws.Range('D36').Value = '0.01804'
$ws.Range('E36').Value = '  +0.14%  '
$ws.Range('D37').Value = '1.228.51'
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').Value = '2.718'
$ws.Range('E38').Value = '  -2.17%  '
$ws.Range('D39').Value = '6.358'
$ws.Range('E39').Value = '  -2.02%  '
$ws.Range('D40').Value = '0.9037'
$ws.Range('E40').Value = '  -0.96%  '
$ws.Range('D41').Value = '0.9991'
$ws.Range('E41').Value = '  -0.15%  '
$ws.Range('E42').Value = '  -0.22%  '
$ws.Range('D43').Value = '65.54'
$ws.Range('E43').Value = '  -1.22%  '
$ws.Range('D44').Value = '7.201'
$ws.Range('E44').Value = '  +0.67%  '
$ws.Range('B45').Value = 'TheSandbox'
$ws.Range('C45').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D45').Value = '0.3995'
$ws.Range('E45').Value = '  -0.49%  '
$ws.Range('B46').Value = 'RenderToken'
$ws.Range('C46').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D46').Value = '1.688'
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('B47').Value = 'BabyDogeCoin'
$ws.Range('C47').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D47').Value = '0.00000000115'
$ws.Range('E47').Value = '  -5.39%  '
$ws.Range('B48').Value = 'EnergySwap'
$ws.Range('C48').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D48').Value = '8.989'
$ws.Range('D49').Value = '0.1142'
$ws.Range('E49').Value = '  +1.50%  '
$ws.Range('D50').Value = '0.05700'
$ws.Range('E50').Value = '  -0.22%  '
$ws.Range('D51').Value = '0.4623'
$ws.Range('E51').Value = '  -0.12%  '

# Restore default (Normal) style so no stray formatting is introduced
$ws.Range("D2:E51").Style = "Normal"
